$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: the old "Url" column (B) becomes the new last column (E);
# Label/Text/Expect all shift one column to the left (B/C/D). ---
$ws.Range("B1").Value = "Label"
$ws.Range("C1").Value = "Text"
$ws.Range("D1").Value = "Expect"
$ws.Range("E1").Value = "Url"

# --- Row 2: GoToUrl target moves from column B to the new column E, and
# now points at seleniumhq instead of w3schools. ---
$ws.Range("B2").ClearContents()
$ws.Range("E2").Value = "http://www.seleniumhq.org/"

# --- Row 3: becomes the new "Fill" step (search box + query). ---
$ws.Range("A3").Value = "Fill"
$ws.Range("B3").Value = "search selenium:"
$ws.Range("C3").Value = "download"

# --- Row 4: Click step now targets the "Go" button. ---
$ws.Range("A4").Value = "Click"
$ws.Range("B4").Value = "Go"
$ws.Range("C4").ClearContents()

# --- Row 5: IsTextVisible now checks for "Selenium IDE". ---
$ws.Range("A5").Value = "IsTextVisible"
$ws.Range("C5").Value = "Selenium IDE"
$ws.Range("D5").Value = $true
$ws.Range("E5").ClearContents()

# --- Row 6: IsTextVisible now checks for "Nonexistent Text". ---
$ws.Range("A6").Value = "IsTextVisible"
$ws.Range("C6").Value = "Nonexistent Text"
$ws.Range("D6").Value = $true
$ws.Range("E6").ClearContents()

# --- Row 7: becomes a Click step on the "Selenium IDE" link. ---
$ws.Range("A7").Value = "Click"
$ws.Range("B7").Value = "Selenium IDE"
$ws.Range("C7").ClearContents()
$ws.Range("D7").ClearContents()

# --- Row 8: Comment moves from column D to column C. ---
$ws.Range("A8").Value = "Comment"
$ws.Range("C8").Value = "Test Finished"
$ws.Range("D8").ClearContents()

# --- Column widths: B/C narrower to fit the new, shorter content. ---
$ws.Columns.Item(2).ColumnWidth = 24.5
$ws.Columns.Item(3).ColumnWidth = 24

# --- Selection moves to B7. ---
$ws.Range("B7").Select() | Out-Null
